# Update column F (dSF) values on Sheet1 to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 4
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = 5
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = 8
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = -6
$ws.Range("F15").Value = -1
$ws.Range("F17").Value = -2
$ws.Range("F18").Value = 5
$ws.Range("F19").Value = -4
$ws.Range("F20").Value = -1
$ws.Range("F24").Value = 2
$ws.Range("F26").Value = 0
$ws.Range("F30").Value = 7
$ws.Range("F32").Value = 5
$ws.Range("F33").Value = -5
$ws.Range("F34").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F38").Value = -2
$ws.Range("F40").Value = 2
$ws.Range("F41").Value = 7
$ws.Range("F42").Value = 2
$ws.Range("F43").Value = 3
